# Append/update timestamp: 2025-11-30 18:31 JST
# Update the "取得日時" (acquired datetime) column (A) for all existing data
# rows in the "ランサーズ" sheet from the previous run timestamp to the new one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldTimestamp = "2025-11-30 18:24:03"
$newTimestamp = "2025-11-30 18:31:40"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
